# Trading update: 2026-02-17 20:15:41
# A new MarketMaking trade (#34) was opened at 20:14:40. The "All Trades"
# history sheet gets a new row for it, and the previously-open trades
# (#29-#33) on that sheet lose their "still open" snapshot fields
# (Capital After / slippage / confidence / entry reason / duration) while
# keeping Status = OPEN and Exit Price = 0. The "MarketMaking" sheet only
# ever tracks the single most recent open trade for that strategy, so its
# old row (trade #29) is replaced by the new trade #34 and the rest of the
# previously-open rows are dropped.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "All Trades" sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

foreach ($r in 30..34) {
    $allTrades.Range("G" + $r).Value() = 0
    $allTrades.Range("K" + $r + ":O" + $r).ClearContents()
    $allTrades.Range("Q" + $r).ClearContents()
}

# New row 35 -> trade #34
$allTrades.Range("A35").Value() = 34

# Force the date to stay plain text (matches the other rows), rather than
# letting Excel auto-convert the literal "2026-02-17" into a date serial.
$allTrades.Range("B35").NumberFormat = "@"
$allTrades.Range("B35").Value() = "2026-02-17"
$allTrades.Range("B35").Style = "Normal"

$allTrades.Range("C35").Value() = "20:14:40"
$allTrades.Range("D35").Value() = "MarketMaking"
$allTrades.Range("E35").Value() = "DOWN"
$allTrades.Range("F35").Value() = 0.95
$allTrades.Range("H35").Value() = "OPEN"
$allTrades.Range("I35").Value() = 0
$allTrades.Range("J35").Value() = 0
$allTrades.Range("K35").Value() = 100
$allTrades.Range("L35").Value() = 0
$allTrades.Range("M35").Value() = 0
$allTrades.Range("N35").Value() = 0.6
$allTrades.Range("O35").Value() = "Normal spread capture: 19600 bps"
$allTrades.Range("Q35").Value() = 0

# ---------------------------------------------------------------------
# "MarketMaking" sheet (open-positions view for this strategy)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Drop the old open trades (#30-#33), keeping just the header + row 2.
$mm.Range("A3:Q6").Delete()

# Row 2 becomes the newly opened trade #34.
$mm.Range("A2").Value() = 34
$mm.Range("C2").Value() = "20:14:40"
$mm.Range("E2").Value() = "DOWN"
$mm.Range("F2").Value() = 0.95
